# Add new vocabulary entries to the "verbanden en formules" glossary sheet,
# fill in a few missing definitions for existing (incomplete) entries, then
# re-sort the whole list alphabetically by column A (term) - matching how a
# user would add rows at the bottom and re-run Data > Sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Promote the column-A style of the four existing entries that are
#    about to gain a definition, from the "incomplete" look (no fill) to
#    the "complete" look (fill), matching the formatting of already-
#    complete rows. Pure formatting - doesn't touch any cell value yet.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)    # Kreukellijn
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial($xlPasteFormats)    # Scheurlijn
$ws.Range("A1").Copy()
$ws.Range("A11").PasteSpecial($xlPasteFormats)   # Tabel
$ws.Range("A1").Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)   # Zaagtand

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)    # Kreukellijn definition style
$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)    # Scheurlijn definition style
$ws.Range("B4").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)   # Zaagtand definition style

# ---------------------------------------------------------------------
# 2) Write the new term names and the new definitions, in the same order
#    the author originally typed them, so newly-created entries line up
#    the same way in the workbook's shared-string table.
# ---------------------------------------------------------------------
$ws.Range("B5").Value2 = "Verticale as, Grafiek, Assenstelsel"    # Kreukellijn
$ws.Range("A21").Value2 = "Globalegrafiek"
$ws.Range("A22").Value2 = "Stijgen"
$ws.Range("A20").Value2 = "Dalen"
$ws.Range("A19").Value2 = "Constant"
$ws.Range("A23").Value2 = "Woordformule"
$ws.Range("B11").Value2 = "Formule, Grafiek, Horizontale as, Verticale as"   # Tabel

# Re-use the already-created "Verticale as, Grafiek, Assenstelsel" string
# for the other two entries that share the exact same definition.
$ws.Range("B9").Value2 = "Verticale as, Grafiek, Assenstelsel"    # Scheurlijn
$ws.Range("B18").Value2 = "Verticale as, Grafiek, Assenstelsel"   # Zaagtand

# ---------------------------------------------------------------------
# 3) Re-sort the whole glossary alphabetically by column A, the way
#    Data > Sort would after adding the rows above. The sort/key ranges
#    mirror the sheet's full historical sort range (A2:B51 / A1:A51),
#    matching the last-used sort remembered on the sheet.
# ---------------------------------------------------------------------
$sortRange = $ws.Range("A2:B51")
$keyRange = $ws.Range("A1:A51")

$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($keyRange)
$sort.SetRange($sortRange)
$sort.Header = 0
$sort.Apply()

# ---------------------------------------------------------------------
# 4) Restore the selected cell the author left active in the file.
# ---------------------------------------------------------------------
$ws.Range("A15").Select()
